$wb = $excel.ActiveWorkbook

# --- ALC sheet updates ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 122.454544
$ws.Range("I5").Value = 122.4
$ws.Range("K5").Value = 122.4
$ws.Range("M5").Value = -7.400000000000006
$ws.Range("H15").Value = 42.06
$ws.Range("I15").Value = 42.06
$ws.Range("K15").Value = 126.18
$ws.Range("M15").Value = 42.81999999999999
$ws.Range("H43").Value = 2107.6924
$ws.Range("I43").Value = 2642.8572
$ws.Range("J43").Value = 1483.3334
$ws.Range("K43").Value = 2642.8572
$ws.Range("L43").Value = 1483.3334
$ws.Range("M43").Value = -2573.8572
$ws.Range("N43").Value = -1621.3334
$ws.Range("H127").Value = 2549.568
$ws.Range("I127").Value = 1279.1
$ws.Range("K127").Value = 3837.3
$ws.Range("M127").Value = 1122.7
$ws.Range("H132").Value = 1588826.6
$ws.Range("I132").Value = 1469.8
$ws.Range("J132").Value = 7694045.5
$ws.Range("K132").Value = 4409.4
$ws.Range("L132").Value = 23082136.5
$ws.Range("M132").Value = -1879.4
$ws.Range("N132").Value = -23087196.5
$ws.Range("H137").Value = 1073.0227
$ws.Range("I137").Value = 923.35895
$ws.Range("J137").Value = 2240.4
$ws.Range("K137").Value = 2770.07685
$ws.Range("L137").Value = 6721.200000000001
$ws.Range("M137").Value = -220.0768500000004
$ws.Range("N137").Value = -11821.2
$ws.Range("H138").Value = 3832.5764
$ws.Range("J138").Value = 4913.673
$ws.Range("L138").Value = 14741.019
$ws.Range("N138").Value = -25021.019

# --- ARM sheet updates ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18452.242
$ws.Range("I32").Value = 14263.024
$ws.Range("J32").Value = 75006.664
$ws.Range("K32").Value = 14263.024
$ws.Range("L32").Value = 75006.664
$ws.Range("M32").Value = -13976.024
$ws.Range("N32").Value = -75580.664
$ws.Range("H63").Value = 2454.7144
$ws.Range("I63").Value = 2132.6785
$ws.Range("J63").Value = 3742.8572
$ws.Range("K63").Value = 2132.6785
$ws.Range("L63").Value = 3742.8572
$ws.Range("M63").Value = -1446.6785
$ws.Range("N63").Value = -5114.8572
$ws.Range("H66").Value = 2454.7144
$ws.Range("I66").Value = 2132.6785
$ws.Range("J66").Value = 3742.8572
$ws.Range("K66").Value = 10663.3925
$ws.Range("L66").Value = 18714.286
$ws.Range("M66").Value = -7231.3925
$ws.Range("N66").Value = -25578.286
$ws.Range("H119").Value = 50000
$ws.Range("J119").Value = 50000
$ws.Range("L119").Value = 50000
$ws.Range("N119").Value = -59676
$ws.Range("H131").Value = 50000
$ws.Range("J131").Value = 50000
$ws.Range("L131").Value = 50000
$ws.Range("N131").Value = -60080

# --- BSM sheet updates ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 41697576
$ws.Range("I20").Value = 47466.46
$ws.Range("J20").Value = 90920430
$ws.Range("K20").Value = 47466.46
$ws.Range("L20").Value = 90920430
$ws.Range("M20").Value = -47219.46
$ws.Range("N20").Value = -90920924
$ws.Range("H86").Value = 7469.3213
$ws.Range("I86").Value = 8569.4
$ws.Range("J86").Value = 4719.125
$ws.Range("K86").Value = 8569.4
$ws.Range("L86").Value = 4719.125
$ws.Range("M86").Value = -7446.4
$ws.Range("N86").Value = -6965.125
$ws.Range("H89").Value = 7469.3213
$ws.Range("I89").Value = 8569.4
$ws.Range("J89").Value = 4719.125
$ws.Range("K89").Value = 42847
$ws.Range("L89").Value = 23595.625
$ws.Range("M89").Value = -37231
$ws.Range("N89").Value = -34827.625
$ws.Range("H122").Value = 49780
$ws.Range("J122").Value = 49780
$ws.Range("L122").Value = 49780
$ws.Range("N122").Value = -59580
$ws.Range("H125").Value = 50000
$ws.Range("J125").Value = 50000
$ws.Range("L125").Value = 50000
$ws.Range("N125").Value = -59840
$ws.Range("H126").Value = 59800
$ws.Range("J126").Value = 59800
$ws.Range("L126").Value = 59800
$ws.Range("N126").Value = -69680

# --- CRP sheet updates ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2786.323
$ws.Range("I31").Value = 2148.578
$ws.Range("J31").Value = 4221.25
$ws.Range("K31").Value = 2148.578
$ws.Range("L31").Value = 4221.25
$ws.Range("M31").Value = -1853.578
$ws.Range("N31").Value = -4811.25
$ws.Range("H34").Value = 2786.323
$ws.Range("I34").Value = 2148.578
$ws.Range("J34").Value = 4221.25
$ws.Range("K34").Value = 2148.578
$ws.Range("L34").Value = 4221.25
$ws.Range("M34").Value = -1946.578
$ws.Range("N34").Value = -4625.25
$ws.Range("H86").Value = 1810.5714
$ws.Range("I86").Value = 1355.3334
$ws.Range("J86").Value = 2417.5557
$ws.Range("K86").Value = 1355.3334
$ws.Range("L86").Value = 2417.5557
$ws.Range("M86").Value = -232.3334
$ws.Range("N86").Value = -4663.5557
$ws.Range("H89").Value = 1810.5714
$ws.Range("I89").Value = 1355.3334
$ws.Range("J89").Value = 2417.5557
$ws.Range("K89").Value = 6776.666999999999
$ws.Range("L89").Value = 12087.7785
$ws.Range("M89").Value = -1160.666999999999
$ws.Range("N89").Value = -23319.7785
$ws.Range("H135").Value = 59333.332
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 59333.332
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 59333.332
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -69473.33199999999

# --- CUL sheet updates ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4694.9375
$ws.Range("I3").Value = 2461.9
$ws.Range("J3").Value = 8416.666999999999
$ws.Range("K3").Value = 7385.700000000001
$ws.Range("L3").Value = 25250.001
$ws.Range("M3").Value = -7273.700000000001
$ws.Range("N3").Value = -25474.001
$ws.Range("H39").Value = 6957.7144
$ws.Range("J39").Value = 6957.7144
$ws.Range("L39").Value = 20873.1432
$ws.Range("N39").Value = -21461.1432
$ws.Range("H74").Value = 4437.5
$ws.Range("J74").Value = 4437.5
$ws.Range("L74").Value = 13312.5
$ws.Range("N74").Value = -15434.5
$ws.Range("H77").Value = 4437.5
$ws.Range("J77").Value = 4437.5
$ws.Range("L77").Value = 39937.5
$ws.Range("N77").Value = -50545.5
$ws.Range("H122").Value = 1451.2162
$ws.Range("J122").Value = 1622.5483
$ws.Range("L122").Value = 14602.9347
$ws.Range("N122").Value = -19502.9347
$ws.Range("H131").Value = 68322.234
$ws.Range("I131").Value = 252857.5
$ws.Range("J131").Value = 39932.19
$ws.Range("K131").Value = 758572.5
$ws.Range("L131").Value = 119796.57
$ws.Range("M131").Value = -753532.5
$ws.Range("N131").Value = -129876.57

# --- GSM sheet updates ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2692.375
$ws.Range("J126").Value = 7733.3335
$ws.Range("L126").Value = 23200.0005
$ws.Range("N126").Value = -28140.0005

# --- LTW sheet updates ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3097.5
$ws.Range("I132").Value = 1971.3871
$ws.Range("J132").Value = 4615.304
$ws.Range("K132").Value = 5914.1613
$ws.Range("L132").Value = 13845.912
$ws.Range("M132").Value = -3384.1613
$ws.Range("N132").Value = -18905.912

# --- WVR sheet updates ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 66650
$ws.Range("J46").Value = 66650
$ws.Range("L46").Value = 66650
$ws.Range("N46").Value = -67112
$ws.Range("H134").Value = 66650
$ws.Range("J134").Value = 66650
$ws.Range("L134").Value = 199950
$ws.Range("N134").Value = -205020
